# Añadida config reflectancias de M2 Drone FB
#
# Adds a new worksheet "M2 Drone FB" (after "Verif Payloads C2") containing
# the reflectance configuration table, and moves the "active sheet" focus
# from "Verif Payloads C2" onto the new sheet.

$wb = $excel.ActiveWorkbook

# --- add the new worksheet right after the last existing sheet ---------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "M2 Drone FB"

# --- fill in the reflectance table --------------------------------------
$newSheet.Range("A1").Value = "Master"
$newSheet.Range("B1").Value = "Blanco"
$newSheet.Range("C1").Value = 1

$newSheet.Range("B2").Value = "Gris"
$newSheet.Range("C2").Value = 0.059658673

$newSheet.Range("B3").Value = "Negro"
$newSheet.Range("C3").Value = 0.048

$newSheet.Range("A4").Value = "Banda 1"
$newSheet.Range("B4").Value = "Blanco"
$newSheet.Range("C4").Value = 1

$newSheet.Range("B5").Value = "Gris"
$newSheet.Range("C5").Value = 0.056465615

$newSheet.Range("B6").Value = "Negro"
$newSheet.Range("C6").Value = 0.049

$newSheet.Range("A7").Value = "Banda 2"
$newSheet.Range("B7").Value = "Blanco"
$newSheet.Range("C7").Value = 1

$newSheet.Range("B8").Value = "Gris"
$newSheet.Range("C8").Value = 0.056221712

$newSheet.Range("B9").Value = "Negro"
$newSheet.Range("C9").Value = 0.044

$newSheet.Range("A10").Value = "Banda 3"
$newSheet.Range("B10").Value = "Blanco"
$newSheet.Range("C10").Value = 1

$newSheet.Range("B11").Value = "Gris"
$newSheet.Range("C11").Value = 0.056885354

$newSheet.Range("B12").Value = "Negro"
$newSheet.Range("C12").Value = 0.047

$newSheet.Range("A13").Value = "Banda 4"
$newSheet.Range("B13").Value = "Blanco"
$newSheet.Range("C13").Value = 1

$newSheet.Range("B14").Value = "Gris"
$newSheet.Range("C14").Value = 0.05918782

$newSheet.Range("B15").Value = "Negro"
$newSheet.Range("C15").Value = 0.046

$newSheet.Range("A16").Value = "Banda 5"
$newSheet.Range("B16").Value = "Blanco"
$newSheet.Range("C16").Value = 1

$newSheet.Range("B17").Value = "Gris"
$newSheet.Range("C17").Value = 0.15

$newSheet.Range("B18").Value = "Negro"
$newSheet.Range("C18").Value = 0.047

# --- previously active sheet ("Verif Payloads C2") loses its selection
#     highlight now that "M2 Drone FB" becomes the active tab -----------
$prevSheet = $wb.Worksheets.Item("Verif Payloads C2")
$null = $prevSheet.Range("A1:C18").Select()

# --- new sheet becomes the active tab / active selection ----------------
$null = $newSheet.Range("D11").Select()
$null = $newSheet.Activate()
